$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values
$ws.Range("B2").Value = 93.823255097523997
$ws.Range("C2").Value = 92.73321619510584
$ws.Range("D2").Value = 92.746566079729874
$ws.Range("E2").Value = 93.834676583843432

# Row 3 data values
$ws.Range("B3").Value = 93.237519419263478
$ws.Range("C3").Value = 93.79004036183251
$ws.Range("D3").Value = 91.340887965378144
$ws.Range("E3").Value = 94.534439327533661

# Update selection to match the new range used in the diff
$ws.Range("B1:E3").Select()
